# ItemDB.xlsx - "Entities" sheet update
# [Jihwan] add Player Animation, add Voxel FX, add Player Attack
#
# The weaponType column (F) for the Weapon rows previously all said
# "Melee"; update each weapon's type to reflect its actual weapon
# category (and introduce the new "TwohandSword" type for the Hammer).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Hammer -> weaponType TwohandSword
$ws.Range("F5").Value = "TwohandSword"
# Row 6 - Sword -> weaponType Sword
$ws.Range("F6").Value = "Sword"
# Row 7 - ChainSaw -> weaponType ChainSaw
$ws.Range("F7").Value = "ChainSaw"
# Row 8 - DarkSword -> weaponType Sword
$ws.Range("F8").Value = "Sword"

# Column F now holds a longer value ("TwohandSword"), so widen it to fit.
$ws.Columns.Item(6).AutoFit() | Out-Null

# Update the selected range to F3:F4 (was H7)
$ws.Range("F3:F4").Select() | Out-Null
